# Add geoposition and coordinates:
#  - Insert a new "Address" column between "Type of Organization" and "Country"
#  - Fill the Address column with the Microsoft Dream Space address for every event row
#  - Fix the "Lattitude" header typo to "Latitude"
#  - Replace the placeholder "NA" Longitude/Latitude values with real coordinates

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# How many data rows currently exist (header row + data rows)
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column at I (pushes Country/Longitude/Lattitude to J/K/L)
$ws.Columns.Item(9).Insert()
$ws.Columns.Item(9).ColumnWidth = 30.8828

# Header row
$ws.Range("I1").Value = "Address"
$ws.Range("L1").Value = "Latitude"

# Data rows: Address, Longitude, Latitude
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "One Microsoft Place, Leopardstown, Dublin, D18 P521"
    $ws.Cells.Item($r, 11).Value = "-6.197155"
    $ws.Cells.Item($r, 12).Value = "53.269037"
}
